$wb = $excel.ActiveWorkbook

# --- Sheet1: "emotion_generation" main data sheet ---
$ws1 = $wb.Worksheets.Item(1)

# Update the (non-shared) formula in I2 to substitute commas with dots for the
# numeric-looking tokens, and to (as in the original author's edit) source the
# emotionIntensity token from column B instead of column C.
$ws1.Range("I2").Formula = '=CONCATENATE("EMOTION(",H2,"){value=[",B2,"]; value_type=[BASICEMOTION]; emotionIntensity=[",SUBSTITUTE(B2,",","."),"]; sourceAggr=[",SUBSTITUTE(D2,",","."),"]; sourceLibid=[",SUBSTITUTE(E2,",","."),"]; sourcePleasure=[",SUBSTITUTE(F2,",","."),"]; sourceUnpleasure=[",SUBSTITUTE(G2,",","."),"]}")'

# Update the shared formula anchored at I3 (covers I3:I32) the same way.
$ws1.Range("I3:I32").Formula = '=CONCATENATE("EMOTION(",H3,"){value=[",B3,"]; value_type=[BASICEMOTION]; emotionIntensity=[",SUBSTITUTE(B3,",","."),"]; sourceAggr=[",SUBSTITUTE(D3,",","."),"]; sourceLibid=[",SUBSTITUTE(E3,",","."),"]; sourcePleasure=[",SUBSTITUTE(F3,",","."),"]; sourceUnpleasure=[",SUBSTITUTE(G3,",","."),"]}")'

# Move the active selection from J2:J32 to I2:I32.
$ws1.Range("I2:I32").Select() | Out-Null

# Apply the page setup (paper size / orientation) that was captured in the
# saved workbook.
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

# --- Sheet2: store the newly generated EMOTION(...) text for reference ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A1").Value = "EMOTION(A06_BEAT_BODO_L01:ANGER){value=[ANGER]; value_type=[BASICEMOTION]; emotionIntensity=[ANGER]; sourceAggr=[0.6]; sourceLibid=[0.2]; sourcePleasure=[0.1]; sourceUnpleasure=[0.6]}"
